# "Initial check-in of translations changes."
#
# The settings sheet stores a row of setting_name/value/comment triples.
# The comment in C1 (next to the "display.title" setting_name/value pair)
# is updated from the stale label "display.title" to the corrected,
# more specific label "display.title.text" used elsewhere in the
# translation tooling.

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("settings")
$wsInitial  = $wb.Worksheets.Item("initial")
$wsSurvey   = $wb.Worksheets.Item("survey")

# Update the stray "display.title" comment/label to "display.title.text".
$wsSettings.Range("C1").Value = "display.title.text"

# Re-create the cursor/selection state left behind by the editing session,
# ending with the "settings" sheet active/selected (matches activeTab=3).
$wsInitial.Activate() | Out-Null
$wsInitial.Range("C2").Select() | Out-Null

$wsSurvey.Activate() | Out-Null
$wsSurvey.Range("F2").Select() | Out-Null

$wsSettings.Activate() | Out-Null
$wsSettings.Range("C2").Select() | Out-Null
